# Add season-record columns (Wins, Losses, Ties) to the player table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): AD1="Wins", AE1="Losses", AF1="Ties" ---
# Copy the formatting (bold, centered, bordered header style) from the
# last existing header cell (AC1) into the new header cells, then set
# their text.
$headerSrc = $ws.Range("AC1")

$headerSrc.Copy($ws.Range("AD1"))
$ws.Range("AD1").Value = "Wins"

$headerSrc.Copy($ws.Range("AE1"))
$ws.Range("AE1").Value = "Losses"

$headerSrc.Copy($ws.Range("AF1"))
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-50): season record for every player's team ---
# Houston's 2014 record: 70 wins, 92 losses, 0 ties.
$wins = 70
$losses = 92
$ties = 0

for ($row = 2; $row -le 50; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins    # column AD
    $ws.Cells.Item($row, 31).Value = $losses  # column AE
    $ws.Cells.Item($row, 32).Value = $ties    # column AF
}
